# edit.ps1 - applies "final touches before submit" edits to
# Graph-Embedding Summary Report.docx via Word COM-interop (iron_native).
#
# Each change below is applied with $d.Content.Find.Execute(...) using a
# narrow, unique anchor string so we do not disturb unrelated text that
# happens to share common words.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "MISS: [$findText]"
    }
    return $ok
}

# 1) "Some of it were in Chinese." -> "All of it in Chinese."
Replace-Text "Some of it were in Chinese." "All of it in Chinese." | Out-Null

# 2) Drop "latter " before "scripts." and change the following sentence's
#    trailing period to a comma ("... the known algorithms." -> "... ,")
Replace-Text "compile two of the latter scripts." "compile two of the scripts." | Out-Null
Replace-Text "the known algorithms." "the known algorithms," | Out-Null

# 3) "Afterward We made" -> "Afterward, We made"
Replace-Text "Afterward We made" "Afterward, We made" | Out-Null

# 4) no textual change to the statistics sentence text itself (only a run
#    split in the original diff) - text remains:
#    " Along with statistical data (such as F1 scores, AUC, etc.) and a confusion matrix."

# 5) no textual change to "project, and get clarifications about the presentation"
#    (only a run split / proofErr wrap in the original diff)

# 6) "we came to the conclusion that the scripts given were not very accurate."
#    -> "we concluded that the scripts given were not giving very accurate results."
#    and "and implementations we found during our research phase" ->
#    "and implementations, during our research phase,"
Replace-Text "and implementations we found during our research phase we came to the conclusion that the scripts given were not very accurate." `
             "and implementations, during our research phase, we concluded that the scripts given were not giving very accurate results." | Out-Null

# 7) Insert a manual line break before "Based on multiple runs of the scripts"
#    and change "the change between the two" -> "the difference between the two"
$r = $d.Content
$ok = $r.Find.Execute("Based on multiple runs of the scripts", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
if ($ok) {
    $brPoint = $d.Range($r.Start, $r.Start)
    $brPoint.InsertBefore([char]11)
} else {
    Write-Output "MISS: [Based on multiple runs of the scripts] (for line break)"
}

Replace-Text "the change between the two" "the difference between the two" | Out-Null

# 8) "walk, because for Node2Vec p=5" -> "walk, due to Node2Vec p=5"
Replace-Text "walk, because for Node2Vec p=5" "walk, due to Node2Vec p=5" | Out-Null

Write-Output "done"
